$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1. Remove the proofErr gramStart/gramEnd markers around
#        "A Cat, a parrot and a bag of seed." (paragraph 4) ---
$pCat = $d.Paragraphs.Item(4)
$pCat.Range.InsertXML("<w:p $wNs><w:r><w:t>A Cat, a parrot and a bag of seed.</w:t></w:r></w:p>")

# --- 2. Define a brand new "A) B) C) ..." (upperLetter) list and apply
#        it to the "A) Man has to..." paragraph (still paragraph 8 here;
#        do this BEFORE splitting the paragraph so Word mints a single,
#        stable numId that we can then reuse explicitly). ---
$pMan = $d.Paragraphs.Item(8)
$pMan.Range.ListFormat.ApplyListTemplateWithLevel($pMan.Range.ListFormat.ListTemplate, $false, 1, $false, 1)
$lt = $pMan.Range.ListFormat.ListTemplate

$styles = @(3,4,2,0,4,2,0,4,2)
for ($i = 1; $i -le 9; $i++) {
    $lvl = $lt.ListLevels.Item($i)
    $lvl.NumberStyle = $styles[$i-1]
    $suffix = if ($i -eq 1) { ")" } else { "." }
    $lvl.NumberFormat = "%" + $i + $suffix
}

# --- 3. Turn that paragraph into three separate list paragraphs sharing
#        the numId that step 2 just minted. ---
$xml  = "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr>"
$xml += "<w:r><w:t>Man has to get to the other side of the river and has room for only one item and himself.</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr>"
$xml += "<w:r><w:t>They are all opposites in nature. The bird fears the cat. The bird loves the seed.</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr>"
$xml += "<w:r><w:t>To safely transport everyone to the other side of the river.</w:t></w:r><w:r><w:t xml:space='preserve'> </w:t></w:r></w:p>"

$pMan.Range.InsertXML($xml)

Write-Output "done"
